$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the third scrum-block header from "Scrumweek  12 - 16 Dec" to "Scrumweek  9 - 13".
# (The old, now-unused shared string drops out and every later shared string's index
# shifts down by one automatically -- which is why C9:C13 end up pointing one lower.)
$ws.Range("A14").Value = "Scrumweek  9 - 13"

# Fill in the "Begin scrum" comments for the new scrum block (column B, rows 15-19).
$ws.Range("B15").Value = "verder met player animations"
$ws.Range("B16").Value = "Verder enemy AI. En wanneer benjamin klaar is met de sprites UI WORK"
$ws.Range("B17").Value = "Werkt lijst met props af"
$ws.Range("B18").Value = "Werkt de level bouw verder af"
$ws.Range("B19").Value = "UI elements"

# Match the wrapped-text row heights these new comments produce.
$ws.Rows.Item(15).RowHeight = 28.8
$ws.Rows.Item(16).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 28.8

# Update the view: scroll position and active selection.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C26").Select()
